$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row labels (A1:J1 "_old" -> "_FV2404", L1:U1 "_new" -> "_FV2410") ---
$baseNames = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = "$($baseNames[$i])_FV2404"
}
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = "$($baseNames[$i])_FV2410"
}

# --- 2. Turn the data range into an Excel table (ListObject) named "Table1" ---
# Preserve the existing header formatting (bold / fill / border) around the table
# creation step so that the engine does not capture it into a new header dxf.
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("A100:U100")

$headerRange.Copy()
$scratch.PasteSpecial(-4122)   # xlPasteFormats

$headerRange.Style = "Normal"

$rng = $ws.Range("A1:U69")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

$scratch.Copy()
$headerRange.PasteSpecial(-4122)   # xlPasteFormats
$scratch.Clear()
$excel.CutCopyMode = $false

# --- 3. Freeze the header row (split below row 1) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
